$d = $word.ActiveDocument

# Namespace-qualified package wrapper used for Range.InsertXML calls.
function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $bodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Splits a paragraph whose full text equals $oldText (with the number-format
# snippet 0" J") into three runs: the unchanged prefix (ending right after
# 0" ), then a run containing "Y", then a run containing the closing quote.
# $rsid is the w:rsidRPr value carried by the original run (kept on the
# first/prefix run, exactly like the source document had it).
function Split-NumberFormatParagraph($paragraph, [string]$prefix, [string]$rsid) {
    $full = $paragraph.Range
    $pStart = $full.Start
    $pEnd = $full.End

    # Range covering the whole paragraph's text, excluding the trailing
    # paragraph mark (End is one past the last visible character).
    $whole = $d.Range($pStart, $pEnd - 1)

    $rsidAttr = ""
    if ($rsid) { $rsidAttr = ' w:rsidRPr="' + $rsid + '"' }

    $run1 = '<w:r' + $rsidAttr + '><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">' + $prefix + '</w:t></w:r>'
    $run2 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Y</w:t></w:r>'
    $run3 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&quot;</w:t></w:r>'

    $xml = New-PkgXml ($run1 + $run2 + $run3)
    $whole.InsertXML($xml)
}

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -eq "The years in columns A and K use the following number format: 0`" J`"`r") {
        Split-NumberFormatParagraph $p "The years in columns A and K use the following number format: 0`" " "00C01BE0"
    }
    elseif ($text -eq "The years in columns A, F and K use the number format 0`" J`"`r") {
        Split-NumberFormatParagraph $p "The years in columns A, F and K use the number format 0`" " "00C01BE0"
    }
}

$word.ActiveDocument.Save()
